# Update the cryptocurrency price/volume table with the latest scraped values.
# Data lives in columns B (Coin), C (Link), D (Price), E (Volume 1h) for rows 2-51.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ref, $value) {
    # Force the cell to stay text (avoids Excel auto-converting numeric-looking
    # strings like "553.91" into floating point numbers).
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $value
}


# Row 2 - Bitcoin
$ws.Range("D2").Value = "57.766.16"
$ws.Range("E2").Value = "  +0.33%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.947.15"
$ws.Range("E3").Value = "  +1.85%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.01%  "

# Row 5 - BNB
Set-TextCell "D5" "553.91"
$ws.Range("E5").Value = "  +0.61%  "

# Row 6 - Solana
Set-TextCell "D6" "133.29"
$ws.Range("E6").Value = "  +10.00%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.04%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +4.11%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.944.07"
$ws.Range("E9").Value = "  +1.90%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +3.81%  "

# Row 11 - Toncoin
Set-TextCell "D11" "4.82"
$ws.Range("E11").Value = "  +2.00%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  +4.04%  "

# Row 13 - ShibaInu
Set-TextCell "D13" "0.0000221"
$ws.Range("E13").Value = "  +4.94%  "

# Row 14 - Avalanche
$ws.Range("E14").Value = "  +4.28%  "

# Row 15 - TRON
$ws.Range("E15").Value = "  +2.86%  "

# Row 16 - WrappedliquidstakedEther2.0
$ws.Range("D16").Value = "3.435.08"
$ws.Range("E16").Value = "  +2.02%  "

# Row 17 - Polkadot
Set-TextCell "D17" "6.96"
$ws.Range("E17").Value = "  +7.32%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.945.46"
$ws.Range("E18").Value = "  +2.17%  "

# Row 19 - WrappedBTC
$ws.Range("D19").Value = "57.723.60"
$ws.Range("E19").Value = "  +0.41%  "

# Row 20 - BitcoinCash
Set-TextCell "D20" "417.34"
$ws.Range("E20").Value = "  +2.27%  "

# Row 21 - Chainlink
Set-TextCell "D21" "13.43"
$ws.Range("E21").Value = "  +5.19%  "

# Row 22 - Polygon
Set-TextCell "D22" "0.703"

# Row 23 - InternetComputer(DFINITY)
Set-TextCell "D23" "13.39"
$ws.Range("E23").Value = "  +6.65%  "

# Row 24 - Uniswap
$ws.Range("E24").Value = "  +4.73%  "

# Row 25 - Litecoin
Set-TextCell "D25" "79.51"
$ws.Range("E25").Value = "  +3.72%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  -0.01%  "

# Row 27 - FirstDigitalUSD
$ws.Range("E27").Value = "  +0.01%  "

# Row 28 - PancakeSwap
Set-TextCell "D28" "2.50"
$ws.Range("E28").Value = "  +1.61%  "

# Row 29 - ImmutableX
$ws.Range("E29").Value = "  +6.93%  "

# Row 30 - RenderToken
$ws.Range("E30").Value = "  +5.91%  "

# Row 31 - EthereumClassic
Set-TextCell "D31" "25.52"
$ws.Range("E31").Value = "  +4.03%  "

# Row 32 - NEARProtocol
Set-TextCell "D32" "5.94"
$ws.Range("E32").Value = "  -1.19%  "

# Row 33 - Hedera
$ws.Range("E33").Value = "  +2.46%  "

# Row 34 - Filecoin (was Mantle)
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell "D34" "5.71"
$ws.Range("E34").Value = "  +7.16%  "

# Row 35 - Mantle (was Filecoin)
$ws.Range("B35").Value = "Mantle"
$ws.Range("C35").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextCell "D35" "0.955"
$ws.Range("E35").Value = "  +6.57%  "

# Row 36 - Stacks
$ws.Range("E36").Value = "  +2.21%  "

# Row 37 - PEPE
$ws.Range("D37").Value = "0.0₃0702"
$ws.Range("E37").Value = "  +14.35%  "

# Row 38 - Cosmos
Set-TextCell "D38" "8.85"
$ws.Range("E38").Value = "  +6.15%  "

# Row 39 - OKB
Set-TextCell "D39" "48.18"
$ws.Range("E39").Value = "  -0.32%  "

# Row 40 - dogwifhat
Set-TextCell "D40" "2.67"
$ws.Range("E40").Value = "  +15.06%  "

# Row 41 - Bittensor
Set-TextCell "D41" "385.30"
$ws.Range("E41").Value = "  +6.53%  "

# Row 42 - Kaspa
$ws.Range("E42").Value = "  +2.72%  "

# Row 43 - VeChain
$ws.Range("E43").Value = "  +1.67%  "

# Row 44 - Maker
$ws.Range("D44").Value = "2.717.22"
$ws.Range("E44").Value = "  +4.41%  "

# Row 45 - USDe
$ws.Range("E45").Value = "  +0.02%  "

# Row 46 - Monero
Set-TextCell "D46" "124.77"
$ws.Range("E46").Value = "  +6.02%  "

# Row 47 - TheGraph
$ws.Range("E47").Value = "  +4.37%  "

# Row 48 - Fetch.AI
$ws.Range("E48").Value = "  +2.93%  "

# Row 49 - Stellar
$ws.Range("E49").Value = "  +1.93%  "

# Row 50 - InjectiveProtocol
Set-TextCell "D50" "22.76"
$ws.Range("E50").Value = "  +2.82%  "

# Row 51 - ThetaToken
$ws.Range("E51").Value = "  +2.92%  "
